$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume snapshot (GitHub Actions scheduled update).
# Price (D) and Volume(1h) (E) columns hold plain text like "145.00" or
# "  -7.19%  "; NumberFormat is forced to text ("@") immediately before each
# assignment so Excel doesn't auto-coerce these into Doubles (which would
# drop trailing zeros / introduce floating point noise).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.965.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.467.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.44%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '547.55'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.34'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.84%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.598'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.465.70'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.32%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -9.23%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.36'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -7.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -7.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.09'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -7.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.903.72'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000163'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -9.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.838.78'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.455.05'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.03'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -7.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.95'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -7.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.18'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -7.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '318.71'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -7.03%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.17'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -6.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.74'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0₃0983'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -8.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.576.39'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.49'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '530.48'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -10.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.30'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -9.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.63'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.149'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.71%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -8.36%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -9.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.85'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -10.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.85'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -10.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.375'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.32'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -7.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '145.22'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -9.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.88'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.28'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -9.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '147.07'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.57'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.85'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -12.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0528'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -10.43%  '

# Rows 50/51 also swap rank order: Mantle now outranks Stellar.
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.583'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -7.24%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0941'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.92%  '
